# Apply crypto price/volume updates as described by the diff.
# All values are written as text (matching the source inlineStr cells),
# using a leading apostrophe to force text interpretation for values that
# look numeric (e.g. "213.22", "0.530"), then resetting the style back to
# Normal so no stray "quote prefix" formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $cell = $ws.Range($range)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "27.708.81"
Set-TextValue "E2" "  +0.99%  "
Set-TextValue "D3" "1.646.62"
Set-TextValue "E3" "  +0.18%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "213.22"
Set-TextValue "E5" "  +0.50%  "
Set-TextValue "D6" "0.530"
Set-TextValue "E6" "  -1.39%  "
Set-TextValue "E7" "  +0.10%  "
Set-TextValue "D8" "23.24"
Set-TextValue "E8" "  +0.49%  "
Set-TextValue "E9" "  +0.67%  "
Set-TextValue "E10" "  +0.59%  "
Set-TextValue "E11" "  +0.40%  "
Set-TextValue "D12" "1.879.28"
Set-TextValue "E12" "  +0.12%  "
Set-TextValue "D13" "1.643.15"
Set-TextValue "E13" "  -0.94%  "
Set-TextValue "D14" "4.04"
Set-TextValue "E14" "  +0.34%  "
Set-TextValue "E15" "  +1.03%  "
Set-TextValue "D16" "64.84"
Set-TextValue "E16" "  +0.79%  "
Set-TextValue "D17" "27.707.91"
Set-TextValue "E17" "  +1.10%  "
Set-TextValue "D18" "232.12"
Set-TextValue "E18" "  +1.59%  "
Set-TextValue "E19" "  +0.86%  "
Set-TextValue "E20" "  +2.01%  "
Set-TextValue "E21" "  +0.09%  "
Set-TextValue "E22" "  -0.58%  "
Set-TextValue "D23" "10.15"
Set-TextValue "E23" "  +8.90%  "
Set-TextValue "E24" "  -2.94%  "
Set-TextValue "D25" "150.06"
Set-TextValue "E25" "  +1.36%  "
Set-TextValue "E26" "  -0.09%  "
Set-TextValue "E27" "  -2.25%  "
Set-TextValue "D28" "15.65"
Set-TextValue "E28" "  +0.74%  "
Set-TextValue "E29" "  +0.02%  "
Set-TextValue "E30" "  +0.49%  "
Set-TextValue "E31" "  +0.06%  "
Set-TextValue "E32" "  +0.94%  "
Set-TextValue "D33" "1.442.68"
Set-TextValue "E33" "  +2.07%  "
Set-TextValue "E34" "  +1.35%  "
Set-TextValue "E35" "  +1.55%  "
Set-TextValue "E36" "  -1.34%  "
Set-TextValue "D37" "0.571"
Set-TextValue "E37" "  +1.50%  "
Set-TextValue "D38" "0.880"
Set-TextValue "E38" "  +0.04%  "
Set-TextValue "D39" "0.0167"
Set-TextValue "E39" "  +0.27%  "
Set-TextValue "D40" "0.884"
Set-TextValue "E40" "  +11.86%  "
Set-TextValue "E41" "  +0.30%  "
Set-TextValue "E42" "  +0.10%  "
Set-TextValue "D43" "67.64"
Set-TextValue "E43" "  +4.63%  "
Set-TextValue "D44" "5.62"
Set-TextValue "E44" "  +2.41%  "
Set-TextValue "E45" "  -0.33%  "
Set-TextValue "E46" "  +1.58%  "
Set-TextValue "D47" "1.788.89"
Set-TextValue "E47" "  +0.00%  "
Set-TextValue "E48" "  +5.02%  "
Set-TextValue "B49" "BabyDogeCoin"
Set-TextValue "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.0₆0107"
Set-TextValue "E49" "  +2.42%  "
Set-TextValue "B50" "Quant"
Set-TextValue "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "85.63"
Set-TextValue "E50" "  -2.16%  "
Set-TextValue "E51" "  +0.15%  "
